# CrewAI Robust Backend Ready!
# Re-label the header row (lowercased / renamed), realign the D:F metric
# columns (shift left + introduce a new "climate change" column) and
# document each column's data type via a legacy cell comment (matches the
# xl/comments/comment1.xml + legacyDrawing VML note pairing added upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data rows: D/E shift down from old E/F, F gets new climate-change values
$ws.Range("D2").Value = 15.65377333333334
$ws.Range("E2").Value = 1732.4853
$ws.Range("F2").Value = 0.00043646874

$ws.Range("D3").Value = 85.60082666666666
$ws.Range("E3").Value = 1129.4168
$ws.Range("F3").Value = 0.0023867781

$ws.Range("D4").Value = 1.24
$ws.Range("E4").Value = 4.72
$ws.Range("F4").Value = 0.00003457449

# --- Column documentation comments (legacy notes, triggers legacyDrawing) --
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
